$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the rich-text shared string "CMS(µ)(2011)" with a plain string
# "CMS(mu)(2011)" for every cell in column L (rows 2-12) that uses it.
$ws.Range("L2:L12").Value = "CMS(mu)(2011)"

# Move the active selection (matches the diff's sheetView selection change).
$ws.Range("L18").Select()
